$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# These cells hold numeric-looking figures stored as text in the workbook.
# Force the "Text" number format first so the new values round-trip as
# shared-string text (matching the source data) rather than being
# reinterpreted as numbers.
$ws.Range("B11:D12").NumberFormat = "@"

# Enterprises density (per 1000 people) - row 11
$ws.Range("B11").Value = "17.56"
$ws.Range("C11").Value = "5.26"
$ws.Range("D11").Value = "22.82"

# Enterprises (% of total) - row 12
$ws.Range("B12").Value = "74.89"
$ws.Range("C12").Value = "22.43"
$ws.Range("D12").Value = "97.31"
